$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue ($ws.Cells.Item(2, 4)) "270.59"

Set-TextValue ($ws.Cells.Item(3, 4)) "21.63"

Set-TextValue ($ws.Cells.Item(4, 4)) "6.327"

Set-TextValue ($ws.Cells.Item(5, 4)) "0.06286"

Set-TextValue ($ws.Cells.Item(6, 4)) "3.565"

Set-TextValue ($ws.Cells.Item(7, 4)) "6.583"

Set-TextValue ($ws.Cells.Item(8, 4)) "1.373"

Set-TextValue ($ws.Cells.Item(9, 4)) "0.8281"

Set-TextValue ($ws.Cells.Item(10, 4)) "0.01376"

Set-TextValue ($ws.Cells.Item(11, 4)) "0.1580"

Set-TextValue ($ws.Cells.Item(12, 4)) "0.08272"

Set-TextValue ($ws.Cells.Item(13, 4)) "0.03383"

Set-TextValue ($ws.Cells.Item(14, 4)) "0.03177"

Set-TextValue ($ws.Cells.Item(15, 2)) "MCDex"
Set-TextValue ($ws.Cells.Item(15, 3)) "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue ($ws.Cells.Item(15, 4)) "4.056"
Set-TextValue ($ws.Cells.Item(15, 5)) "14MCDexMCB"

Set-TextValue ($ws.Cells.Item(16, 2)) "BitMartToken"
Set-TextValue ($ws.Cells.Item(16, 3)) "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue ($ws.Cells.Item(16, 4)) "0.09402"
Set-TextValue ($ws.Cells.Item(16, 5)) "15BitMartTokenBMX"

Set-TextValue ($ws.Cells.Item(17, 4)) "0.001639"

Set-TextValue ($ws.Cells.Item(18, 4)) "0.04679"

Set-TextValue ($ws.Cells.Item(19, 4)) "0.006286"

Set-TextValue ($ws.Cells.Item(20, 4)) "0.005943"

Set-TextValue ($ws.Cells.Item(21, 4)) "0.001062"

Set-TextValue ($ws.Cells.Item(22, 4)) "0.0001491"

Set-TextValue ($ws.Cells.Item(23, 4)) "3.752"

Set-TextValue ($ws.Cells.Item(25, 4)) "0.3301"

Set-TextValue ($ws.Cells.Item(26, 4)) "0.1249"

Set-TextValue ($ws.Cells.Item(28, 4)) "0.0002712"

Set-TextValue ($ws.Cells.Item(40, 4)) "0.04707"

Set-TextValue ($ws.Cells.Item(41, 4)) "0.007097"

Set-TextValue ($ws.Cells.Item(42, 4)) "0.1176"

Set-TextValue ($ws.Cells.Item(43, 4)) "0.003638"

Set-TextValue ($ws.Cells.Item(44, 4)) "0.01168"

Set-TextValue ($ws.Cells.Item(45, 4)) "0.00005979"

Set-TextValue ($ws.Cells.Item(46, 4)) "0.0009829"

Set-TextValue ($ws.Cells.Item(47, 4)) "0.00000000746"

Set-TextValue ($ws.Cells.Item(48, 4)) "0.7778"

Set-TextValue ($ws.Cells.Item(49, 4)) "0.002394"

Set-TextValue ($ws.Cells.Item(50, 4)) "0.00002386"

Set-TextValue ($ws.Cells.Item(51, 4)) "0.01233"

Write-Host "Updated crypto price/symbol list."
